$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add header for new year column AH (2023), matching style/formatting of the
# existing header cells (e.g. AG1: bold, bordered, centered/top-aligned).
$ws.Range("AG1").Copy()
$ws.Range("AH1").PasteSpecial(-4122)
$ws.Range("AH1").Value = "'2023"

# Add the two new data values for 2023 (plain numbers, like the other years)
$ws.Range("AH2").Value = 97604.7
$ws.Range("AH3").Value = 8621.4
